$d = $word.ActiveDocument

# Locate the paragraph holding the phone number; the new "Github" line is
# inserted directly after it.
$phonePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("+44 7488 495712")) {
        $phonePara = $p
        break
    }
}

# Collapse to a zero-length range just before the paragraph mark of the
# phone-number paragraph so the new content is inserted as a new
# paragraph immediately after it (without disturbing the following
# empty paragraph at the end of the document).
$insertPoint = $phonePara.Range.End - 1
$insertRange = $d.Range($insertPoint, $insertPoint)

# Build the new paragraph as a WordprocessingML fragment so we get the
# same run/proofErr layout Word itself produces when it spell-checks the
# non-dictionary word "Github" (a spellStart/spellEnd pair wrapping its
# own run, followed by separate runs for the separator and the URL).
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="p3"/>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Github</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> : </w:t>
            </w:r>
            <w:r>
              <w:t>https://github.com/vishals18/oxford-risk-internship/tree/main</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$insertRange.InsertXML($xml)
